$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Metadata schema ID" / "Header schema ID" field which lives in row 1
# (A1 = "Header schema ID", B1 = the schema UUID). Deleting the whole row shifts
# everything else up by one, matching the target layout.
$ws.Rows("1:1").Delete()
